# Append new scraped job listing (2026-01-02 18:27 JST run) and shift
# existing rows down to keep the sheet ordered by priority score.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-01-02 18:27:29"

# Final desired state for data rows 2..9 (row 1 is the header, untouched).
$rows = @(
    @{ B = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5460562"; G = 435; H = "🔥AI,Ai ◆ツール,開発" },
    @{ B = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5423720"; G = 385; H = "🔥AI,Ai ◆効率化" },
    @{ B = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5460563"; G = 220; H = "◆開発,システム開発 ◇管理" },
    @{ B = "【介護業務効率化】研修事業の自動化を実現するプロ募集"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5464016"; G = 153; H = "◆効率化,自動化" },
    @{ B = "ホットペッパービューティーブログ一括投稿システム開発"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5455160"; G = 113; H = "◆開発,システム開発" },
    @{ B = "複数WEBサイトへの日記一括投稿ツールの修正 or 新規作成をお願いしたいです"; C = "システム開発"; D = "1,000 ~ 5,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5463948"; G = 90; H = "◆ツール ◇サイト" },
    @{ B = "【報告書自動化】GASで効率的な作成フローを実現!"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5464025"; G = 80; H = "◆自動化" },
    @{ B = "進行管理およびチームディレクションを担当"; C = "システム開発"; D = "~ 5,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5418064"; G = 30; H = "◇管理" }
)

# Clear any existing hyperlinks up front; they will be rebuilt from scratch
# once every cell holds its final value, so ranges/targets stay consistent.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# Re-create the hyperlinks for column F (rows 2..9) in order, and restore
# the "Hyperlink" cell style that Value assignment alone does not set.
$r = 2
foreach ($row in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $row.F)
    $cell.Style = "Hyperlink"
    $r = $r + 1
}
